$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39, -49.8303533733402,  -1.285958296698702,  0.03012870743856047, 0.002715833562876155, 0.03149110810101963),
    @(40, -44.45422979399342, -0.9429480844145793, 0.03727361920264568, 0.004112031532115542, 0.0389895444232617),
    @(41, -49.12644518608512, -1.234518913264335,  0.08183332923410236, 0.002938813338652266, 0.03307793574342761),
    @(42, -49.08736175233031, -1.262833530135628,  0.03039696420631934, 0.002849600249672071, 0.03176445169193558),
    @(43, -48.72457008755512, -1.204407234879717,  0.03203233451463891, 0.003057073877455373, 0.03426768450626295),
    @(44, -42.10659406711788, -0.8375361046424108, 0.04028077999362526, 0.004487459974258948, 0.04242369853510088),
    @(45, -50.74085032341767, -1.35207691598998,   0.05923577998107189, 0.002453689080871777, 0.03743885377620029),
    @(46, -43.85653625902743, -0.9193241304295263, 0.03756831901965291, 0.004167149281924271, 0.03937916771150821)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# Row 9 is a brand-new row: give column A the same style as the rest of
# the A-column data cells (bold, bordered, centered) by copying A8's format.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

